$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (rows 3-8): replace static "NBL017x" product codes with a
# generated "PK_SAC_#####" formula ---------------------------------------
# Clear existing formatting on B3:B8 and re-apply a plain thin border
# (matches the border already used throughout the sheet) without carrying
# over the old font / alignment.
$rngB = $ws.Range("B3:B8")
$rngB.Style = "Normal"
$rngB.Borders.Color = 0
$rngB.Borders.Weight = 2
$rngB.Borders.LineStyle = 1

# B3 holds its own (non-shared) formula, B4:B8 share one formula.
$ws.Range("B3").Formula = '="PK_SAC_"&TEXT(RANDBETWEEN(0,99999),"00000")'
$ws.Range("B4:B8").Formula = '="PK_SAC_"&TEXT(RANDBETWEEN(0,99999),"00000")'

# --- Row heights: rows 3-8 go back to the sheet's default height --------
$ws.Rows("3:8").AutoFit()

# --- Column widths --------------------------------------------------------
$ws.Columns("A").ColumnWidth = 37
$ws.Columns("R").ColumnWidth = 18.1666666666666667

# --- Selection -------------------------------------------------------------
$ws.Range("B26").Select()
